$d = $word.ActiveDocument

function Split-AtBreak($findText, $replaceText) {
    $range = $d.Content
    $result = $range.Find.Execute($findText, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, $replaceText, 2)
    if (-not $result) {
        throw "Find/Replace failed for: $findText"
    }
}

# 1) Objetivos paragraph: split into three runs/text-lines with manual line breaks.
Split-AtBreak "tabela periódica.-Capacitar" "tabela periódica.^l-Capacitar"
Split-AtBreak "aplicações- capacitar" "aplicações^l- capacitar"

# 2) Programa resumido paragraph: split into two lines.
Split-AtBreak "Grupo 13.- Metais" "Grupo 13.^l- Metais"

# 3) Programa paragraph: split into two lines.
Split-AtBreak "Complexos.Relacionar" "Complexos.^lRelacionar"

# 4) Avaliação - Método run: split into two lines.
Split-AtBreak "objetivadas.Duas" "objetivadas.^lDuas"
